$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "BALANCE : 10000524.0"

$newRows = @(
    @("2025-10-21 20:45:54", "Withdraw", 1, 399),
    @("2025-10-21 20:46:47", "Deposit", 2, 401),
    @("2025-10-21 20:53:11", "Deposit", 123, 524),
    @("2025-10-21 21:25:42", "Deposit", 10000000, 10000524)
)

$r = 32
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    $rowRange = $ws.Range("A${r}:D${r}")
    $rowRange.NumberFormat = "General"
    $r++
}
